# Update "想去人数" (interest count) values in two worksheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1102
$wsExhibition.Range("F4").Value = 2528
$wsExhibition.Range("F5").Value = 216

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1102
$wsAll.Range("F6").Value = 2528
$wsAll.Range("F8").Value = 216
